# Auto-generated edit script: updates crypto price (D) and volume% (E) columns
# to match the scraped values from the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.364.97"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.566.91"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.790.20"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.566.53"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "27.388.20"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "0.0₃0688"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").Value = "1.371.05"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.968"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.531"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("D47").Value = "1.701.83"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").Value = "0.0₇0991"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("E51").Value = "  -0.87%  "
